# This workbook has two sheets ("展览" and "全部类型") that carry an
# identical table of convention-listing rows. The scraped data moved on:
# the two oldest events (rows 2 and 3) dropped off the list, every other
# event shifted up two rows, the trailing two rows disappeared, a handful
# of "want-to-go" counters (column F) ticked up, and the ticket-price cell
# for the event that is now row 2 flipped from a price to "不可售".

$wb = $excel.ActiveWorkbook

# New value for column F (want-to-go count), keyed by the *new* row number
# (i.e. after the first two old rows have been dropped and everything else
# shifted up by two).
$fUpdates = @{
    4  = 1768
    6  = 757
    9  = 33
    14 = 124
    15 = 142
    16 = 4135
    17 = 12
    19 = 458
    20 = 394
    21 = 960
    22 = 1218
    24 = 31
    25 = 39
    26 = 1885
    27 = 59
    28 = 46
    29 = 81
    31 = 24
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Drop the two oldest events; everything below shifts up two rows and
    # the sheet's used range shrinks from 33 to 31 rows automatically.
    $ws.Range("A2:I3").EntireRow.Delete() | Out-Null

    # Column A is just a running index (row number - 1); restore it since
    # the row delete shifted those values along with everything else.
    for ($r = 1; $r -le 31; $r++) {
        $ws.Cells.Item($r, 1).Value2 = $r - 1
    }

    # Apply the updated "want-to-go" counts.
    foreach ($r in $fUpdates.Keys) {
        $ws.Cells.Item($r, 6).Value2 = $fUpdates[$r]
    }

    # The event now sitting in row 2 is no longer sellable; its minimum
    # price cell becomes the text "不可售" instead of a numeric price.
    $ws.Cells.Item(2, 7).Value = "不可售"
}
